$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.717.41"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "3.755.71"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'602.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'169.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("D7").Value = "3.753.76"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("D11").Value = "'6.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "'38.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "4.384.08"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "3.753.04"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "68.750.46"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").Value = "'7.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'17.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").Value = "'497.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").Value = "'10.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.39%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'85.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("D26").Value = "'2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "'12.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("D28").Value = "'10.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +5.24%  "
$ws.Range("D31").Value = "'2.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").Value = "'7.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").Value = "'32.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "3.903.75"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "3.690.03"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "'0.108"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'5.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.325"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").Value = "'444.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").Value = "'48.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "'2.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D48").Value = "'40.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").Value = "2.847.62"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "'142.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("E51").Value = "  +2.96%  "